$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "28.252.63"
Set-TextValue $ws.Range("E2") "  -5.32%  "
Set-TextValue $ws.Range("D3") "1.838.34"
Set-TextValue $ws.Range("E3") "  -5.01%  "
Set-TextValue $ws.Range("D4") "1.003"
Set-TextValue $ws.Range("E4") "  -0.38%  "
Set-TextValue $ws.Range("D5") "330.10"
Set-TextValue $ws.Range("E5") "  -1.83%  "
Set-TextValue $ws.Range("D6") "1.002"
Set-TextValue $ws.Range("E6") "  -0.47%  "
Set-TextValue $ws.Range("D7") "0.4596"
Set-TextValue $ws.Range("E7") "  -4.88%  "
Set-TextValue $ws.Range("D8") "0.3862"
Set-TextValue $ws.Range("E8") "  -5.94%  "
Set-TextValue $ws.Range("B9") "Dogecoin"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D9") "0.07855"
Set-TextValue $ws.Range("E9") "  -3.85%  "
Set-TextValue $ws.Range("B10") "Polygon"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D10") "0.9612"
Set-TextValue $ws.Range("E10") "  -5.16%  "
Set-TextValue $ws.Range("B11") "Solana"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws.Range("D11") "21.94"
Set-TextValue $ws.Range("E11") "  -7.25%  "
Set-TextValue $ws.Range("B12") "WrappedEther"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D12") "1.805.75"
Set-TextValue $ws.Range("E12") "  -7.22%  "
Set-TextValue $ws.Range("D13") "5.714"
Set-TextValue $ws.Range("E13") "  -5.97%  "
Set-TextValue $ws.Range("B14") "Chainlink"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D14") "6.917"
Set-TextValue $ws.Range("E14") "  -4.91%  "
Set-TextValue $ws.Range("B15") "TRON"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D15") "0.06852"
Set-TextValue $ws.Range("E15") "  +0.09%  "
Set-TextValue $ws.Range("B16") "BinanceUSD"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D16") "1.002"
Set-TextValue $ws.Range("E16") "  -0.53%  "
Set-TextValue $ws.Range("B17") "Litecoin"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D17") "86.90"
Set-TextValue $ws.Range("E17") "  -4.40%  "
Set-TextValue $ws.Range("B18") "ShibaInu"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D18") "0.000009936"
Set-TextValue $ws.Range("E18") "  -3.86%  "
Set-TextValue $ws.Range("B19") "Avalanche"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D19") "16.90"
Set-TextValue $ws.Range("E19") "  -4.91%  "
Set-TextValue $ws.Range("B20") "Dai"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D20") "1.002"
Set-TextValue $ws.Range("E20") "  -0.38%  "
Set-TextValue $ws.Range("B21") "WrappedBTC"
Set-TextValue $ws.Range("C21") "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D21") "28.284.57"
Set-TextValue $ws.Range("E21") "  -5.14%  "
Set-TextValue $ws.Range("B22") "Uniswap"
Set-TextValue $ws.Range("C22") "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D22") "5.345"
Set-TextValue $ws.Range("E22") "  -4.91%  "
Set-TextValue $ws.Range("B23") "Cosmos"
Set-TextValue $ws.Range("C23") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D23") "10.96"
Set-TextValue $ws.Range("E23") "  -7.56%  "
Set-TextValue $ws.Range("B24") "Toncoin"
Set-TextValue $ws.Range("C24") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D24") "2.142"
Set-TextValue $ws.Range("E24") "  -1.52%  "
Set-TextValue $ws.Range("B25") "WrappedliquidstakedEther2.0"
Set-TextValue $ws.Range("C25") "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D25") "2.091.02"
Set-TextValue $ws.Range("E25") "  -3.44%  "
Set-TextValue $ws.Range("B26") "Monero"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D26") "153.36"
Set-TextValue $ws.Range("E26") "  -2.06%  "
Set-TextValue $ws.Range("B27") "EthereumClassic"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D27") "19.22"
Set-TextValue $ws.Range("E27") "  -4.15%  "
Set-TextValue $ws.Range("B28") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D28") "5.733"
Set-TextValue $ws.Range("E28") "  -13.21%  "
Set-TextValue $ws.Range("B29") "LidoDAOToken"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D29") "1.983"
Set-TextValue $ws.Range("E29") "  -5.30%  "
Set-TextValue $ws.Range("B30") "BitcoinCash"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D30") "116.98"
Set-TextValue $ws.Range("E30") "  -3.29%  "
Set-TextValue $ws.Range("B31") "ImmutableX"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D31") "0.9427"
Set-TextValue $ws.Range("E31") "  -6.36%  "
Set-TextValue $ws.Range("B32") "Stellar"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D32") "0.09296"
Set-TextValue $ws.Range("E32") "  -3.35%  "
Set-TextValue $ws.Range("B33") "Filecoin"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D33") "5.283"
Set-TextValue $ws.Range("E33") "  -4.79%  "
Set-TextValue $ws.Range("B34") "HuobiToken"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D34") "3.446"
Set-TextValue $ws.Range("E34") "  -2.39%  "
Set-TextValue $ws.Range("B35") "ARBITRUM"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D35") "1.326"
Set-TextValue $ws.Range("E35") "  -6.26%  "
Set-TextValue $ws.Range("B36") "Hedera"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D36") "0.06006"
Set-TextValue $ws.Range("E36") "  -8.56%  "
Set-TextValue $ws.Range("B37") "VeChain"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D37") "0.02150"
Set-TextValue $ws.Range("E37") "  -5.69%  "
Set-TextValue $ws.Range("B38") "TrustWalletToken"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D38") "1.146"
Set-TextValue $ws.Range("E38") "  -4.74%  "
Set-TextValue $ws.Range("B39") "Frax"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D39") "1.001"
Set-TextValue $ws.Range("E39") "  -0.46%  "
Set-TextValue $ws.Range("B40") "FraxShare"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D40") "7.628"
Set-TextValue $ws.Range("E40") "  -3.82%  "
Set-TextValue $ws.Range("B41") "TheSandbox"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D41") "0.5620"
Set-TextValue $ws.Range("E41") "  -5.70%  "
Set-TextValue $ws.Range("B42") "Aptos"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D42") "9.988"
Set-TextValue $ws.Range("E42") "  -6.95%  "
Set-TextValue $ws.Range("B43") "Algorand"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D43") "0.1778"
Set-TextValue $ws.Range("E43") "  -3.66%  "
Set-TextValue $ws.Range("B44") "WEMIXToken"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D44") "1.252"
Set-TextValue $ws.Range("E44") "  -1.75%  "
Set-TextValue $ws.Range("B45") "RenderToken"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D45") "2.270"
Set-TextValue $ws.Range("E45") "  -8.28%  "
Set-TextValue $ws.Range("B46") "EnergySwap"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D46") "11.64"
Set-TextValue $ws.Range("E46") "  -5.23%  "
Set-TextValue $ws.Range("B47") "Decentraland"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D47") "0.5288"
Set-TextValue $ws.Range("E47") "  -4.81%  "
Set-TextValue $ws.Range("B48") "Cronos"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D48") "0.07030"
Set-TextValue $ws.Range("E48") "  -5.94%  "
Set-TextValue $ws.Range("B49") "NEARProtocol"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D49") "1.833"
Set-TextValue $ws.Range("E49") "  -7.73%  "
Set-TextValue $ws.Range("B50") "Quant"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D50") "112.48"
Set-TextValue $ws.Range("E50") "  -3.95%  "
Set-TextValue $ws.Range("B51") "PaxDollar"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D51") "1.001"
Set-TextValue $ws.Range("E51") "  -0.65%  "
